# Update date line
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-09 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-10 Sunday", 2)

# Table 1 cell-by-cell replacements (row,col are 1-based; only data rows are 1,5,9,13,17)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "30÷2=15, 0"
$t.Cell(1,2).Range.Text = "91÷8=11, 3"
$t.Cell(1,3).Range.Text = "47÷8=5, 7"
$t.Cell(1,4).Range.Text = "91÷6=15, 1"
$t.Cell(1,5).Range.Text = "45÷8=5, 5"

$t.Cell(5,1).Range.Text = "87÷4=21, 3"
$t.Cell(5,2).Range.Text = "22÷9=2, 4"
$t.Cell(5,3).Range.Text = "82÷4=20, 2"
$t.Cell(5,4).Range.Text = "68÷2=34, 0"
$t.Cell(5,5).Range.Text = "12÷5=2, 2"

$t.Cell(9,1).Range.Text = "88÷5=17, 3"
$t.Cell(9,2).Range.Text = "24÷3=8, 0"
$t.Cell(9,3).Range.Text = "63÷9=7, 0"
$t.Cell(9,4).Range.Text = "51÷8=6, 3"
$t.Cell(9,5).Range.Text = "50÷4=12, 2"

$t.Cell(13,1).Range.Text = "66÷9=7, 3"
$t.Cell(13,2).Range.Text = "90÷6=15, 0"
$t.Cell(13,3).Range.Text = "50÷8=6, 2"
$t.Cell(13,4).Range.Text = "93÷2=46, 1"
$t.Cell(13,5).Range.Text = "96÷9=10, 6"

$t.Cell(17,1).Range.Text = "76÷7=10, 6"
$t.Cell(17,2).Range.Text = "91÷2=45, 1"
$t.Cell(17,3).Range.Text = "99÷6=16, 3"
$t.Cell(17,4).Range.Text = "54÷3=18, 0"
$t.Cell(17,5).Range.Text = "48÷4=12, 0"
